# Journal de bord CPNVoiturage - mise a jour de la doc
# Fill in rows 50-59 of the "Page 1" sheet (table Tableau1) with the
# new journal entries, and move the active selection to F60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 50
$ws.Range("C50").Value = "Affichage dynamique des données sur le profil"
$ws.Range("D50").Value = 44270
$ws.Range("E50").Value = 60

# Row 51
$ws.Range("C51").Value = "Requetes et création du tableau pour l'affichage du profil"
$ws.Range("D51").Value = 44270
$ws.Range("E51").Value = 60

# Row 52
$ws.Range("C52").Value = "Utilisation des Ids des périodes pour le profil"
$ws.Range("D52").Value = 44273
$ws.Range("E52").Value = 45

# Row 53
$ws.Range("C53").Value = "Gestion des colonnes nécessaires de csv de EDT"
$ws.Range("D53").Value = 44273
$ws.Range("E53").Value = 45

# Row 54
$ws.Range("C54").Value = "Affichage de la page admin"
$ws.Range("D54").Value = 44273
$ws.Range("E54").Value = 30

# Row 55
$ws.Range("C55").Value = "Récuperation et processing du csv"
$ws.Range("D55").Value = 44273
$ws.Range("E55").Value = 80

# Row 56
$ws.Range("C56").Value = "Récuperation pour préparer la requête de EDT"
$ws.Range("D56").Value = 44273
$ws.Range("E56").Value = 45

# Row 57
$ws.Range("C57").Value = "Fin de l'importation de EDT"
$ws.Range("D57").Value = 44274
$ws.Range("E57").Value = 40

# Row 58
$ws.Range("C58").Value = "Corrections sur le login"
$ws.Range("D58").Value = 44274
$ws.Range("E58").Value = 15
$ws.Range("F58").Value = "3 caractères max (coté client et serveur) pour l'acronyme et ajout du signout"

# Row 59
$ws.Range("C59").Value = "Refonte d'une requête"
$ws.Range("D59").Value = 44274
$ws.Range("E59").Value = 70
$ws.Range("F59").Value = "Refonte de la requête pour afficher les voitures car il manquait certaines personnes seules"

# Update the active selection to reflect the latest edited cell
$ws.Range("F60").Select()
